$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add the new "value_unit" column (K) with its header + numeric data.
# ---------------------------------------------------------------------------
$ws.Range("K1").Value = "value_unit"

# Give the new header cell the same bold / filled look as the other header
# cells (copy the format from J1) and then drop the border so it matches
# the "no border" header style used for this new trailing column.
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Borders.Item(1).LineStyle = 0
$ws.Range("K1").Borders.Item(2).LineStyle = 0
$ws.Range("K1").Borders.Item(3).LineStyle = 0
$ws.Range("K1").Borders.Item(4).LineStyle = 0
$ws.Range("K1").Borders.Item(5).LineStyle = 0
$ws.Range("K1").Borders.Item(6).LineStyle = 0

$kValues = @{
    2  = 419430
    3  = 420780
    4  = 421980
    5  = 189280
    6  = 190350
    7  = 190700
    8  = 208920
    9  = 210200
    10 = 210300
    11 = 458080
    12 = 460860
    13 = 462340
    14 = 41780
    15 = 41850
    16 = 41970
    17 = 195270
    18 = 199500
    19 = 199720
    20 = 37500
    21 = 37570
    22 = 37610
    23 = 174430
    24 = 174530
    25 = 175830
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 11).Value = $kValues[$row]
}

# ---------------------------------------------------------------------------
# Drop the stray formatting previously carried on these cells (it was a
# no-op fill/border combination visually identical to the default style).
# ---------------------------------------------------------------------------
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("G20").ClearFormats()
$ws.Range("G24").ClearFormats()

# ---------------------------------------------------------------------------
# Update the active selection to match the saved view.
# ---------------------------------------------------------------------------
$ws.Range("N17").Select() | Out-Null
